# New PO forecast model
# Updates three worksheets:
#   - "Weekly Quantity": append one new weekly data point (row 62)
#   - "Monthly Trend":   append one new monthly data point (row 21)
#   - "PO Forecast":     recompute the cumulative forecast curve - every
#                        existing quantity (col B) changes, and the tail of
#                        the date series (col A, rows 62-69) shifts forward
#                        in time with one brand-new row (70) appended.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Weekly Quantity": new row 62
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Cells.Item(62,1).Value = 45662.99999999999
$wsWeekly.Cells.Item(62,2).Value = 150
$wsWeekly.Cells.Item(62,1).NumberFormat = $wsWeekly.Cells.Item(61,1).NumberFormat

# ---------------------------------------------------------------------
# Sheet "Monthly Trend": new row 21
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Cells.Item(21,1).Value = 45688.99999999999
$wsMonthly.Cells.Item(21,2).Value = 150
$wsMonthly.Cells.Item(21,1).NumberFormat = $wsMonthly.Cells.Item(20,1).NumberFormat

# ---------------------------------------------------------------------
# Sheet "PO Forecast": update B column (rows 2-61) and A/B columns (rows 62-70)
# ---------------------------------------------------------------------
$wsPO = $wb.Worksheets.Item("PO Forecast")

# Rows 2-61: only the Qty (column B) changes
$wsPO.Cells.Item(2,2).Value = 158
$wsPO.Cells.Item(3,2).Value = 164
$wsPO.Cells.Item(4,2).Value = 174
$wsPO.Cells.Item(5,2).Value = 180
$wsPO.Cells.Item(6,2).Value = 185
$wsPO.Cells.Item(7,2).Value = 191
$wsPO.Cells.Item(8,2).Value = 196
$wsPO.Cells.Item(9,2).Value = 201
$wsPO.Cells.Item(10,2).Value = 207
$wsPO.Cells.Item(11,2).Value = 212
$wsPO.Cells.Item(12,2).Value = 217
$wsPO.Cells.Item(13,2).Value = 223
$wsPO.Cells.Item(14,2).Value = 228
$wsPO.Cells.Item(15,2).Value = 239
$wsPO.Cells.Item(16,2).Value = 244
$wsPO.Cells.Item(17,2).Value = 250
$wsPO.Cells.Item(18,2).Value = 255
$wsPO.Cells.Item(19,2).Value = 261
$wsPO.Cells.Item(20,2).Value = 277
$wsPO.Cells.Item(21,2).Value = 282
$wsPO.Cells.Item(22,2).Value = 288
$wsPO.Cells.Item(23,2).Value = 293
$wsPO.Cells.Item(24,2).Value = 304
$wsPO.Cells.Item(25,2).Value = 325
$wsPO.Cells.Item(26,2).Value = 331
$wsPO.Cells.Item(27,2).Value = 336
$wsPO.Cells.Item(28,2).Value = 341
$wsPO.Cells.Item(29,2).Value = 347
$wsPO.Cells.Item(30,2).Value = 352
$wsPO.Cells.Item(31,2).Value = 358
$wsPO.Cells.Item(32,2).Value = 368
$wsPO.Cells.Item(33,2).Value = 374
$wsPO.Cells.Item(34,2).Value = 379
$wsPO.Cells.Item(35,2).Value = 384
$wsPO.Cells.Item(36,2).Value = 465
$wsPO.Cells.Item(37,2).Value = 476
$wsPO.Cells.Item(38,2).Value = 481
$wsPO.Cells.Item(39,2).Value = 487
$wsPO.Cells.Item(40,2).Value = 492
$wsPO.Cells.Item(41,2).Value = 498
$wsPO.Cells.Item(42,2).Value = 503
$wsPO.Cells.Item(43,2).Value = 508
$wsPO.Cells.Item(44,2).Value = 514
$wsPO.Cells.Item(45,2).Value = 519
$wsPO.Cells.Item(46,2).Value = 530
$wsPO.Cells.Item(47,2).Value = 535
$wsPO.Cells.Item(48,2).Value = 541
$wsPO.Cells.Item(49,2).Value = 546
$wsPO.Cells.Item(50,2).Value = 552
$wsPO.Cells.Item(51,2).Value = 562
$wsPO.Cells.Item(52,2).Value = 568
$wsPO.Cells.Item(53,2).Value = 573
$wsPO.Cells.Item(54,2).Value = 578
$wsPO.Cells.Item(55,2).Value = 584
$wsPO.Cells.Item(56,2).Value = 616
$wsPO.Cells.Item(57,2).Value = 622
$wsPO.Cells.Item(58,2).Value = 654
$wsPO.Cells.Item(59,2).Value = 670
$wsPO.Cells.Item(60,2).Value = 675
$wsPO.Cells.Item(61,2).Value = 681

# Rows 62-69: existing rows whose date (A) and qty (B) both shift
$wsPO.Cells.Item(62,1).Value = 45662.99999999999
$wsPO.Cells.Item(62,2).Value = 719
$wsPO.Cells.Item(63,1).Value = 45669.99999999999
$wsPO.Cells.Item(63,2).Value = 724
$wsPO.Cells.Item(64,1).Value = 45676.99999999999
$wsPO.Cells.Item(64,2).Value = 729
$wsPO.Cells.Item(65,1).Value = 45683.99999999999
$wsPO.Cells.Item(65,2).Value = 735
$wsPO.Cells.Item(66,1).Value = 45690.99999999999
$wsPO.Cells.Item(66,2).Value = 740
$wsPO.Cells.Item(67,1).Value = 45697.99999999999
$wsPO.Cells.Item(67,2).Value = 746
$wsPO.Cells.Item(68,1).Value = 45704.99999999999
$wsPO.Cells.Item(68,2).Value = 751
$wsPO.Cells.Item(69,1).Value = 45711.99999999999
$wsPO.Cells.Item(69,2).Value = 756

# Row 70: brand-new row, copy the date number format from row 69
$wsPO.Cells.Item(70,1).Value = 45718.99999999999
$wsPO.Cells.Item(70,2).Value = 762
$wsPO.Cells.Item(70,1).NumberFormat = $wsPO.Cells.Item(69,1).NumberFormat
